$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Step2 test case (row 20-21): add the extra b[] element columns (J, K)
# that addAll(a[], b, c[]) now produces at runtime.
$ws.Range("J20").Value = '_res_.$Step2[3]'
$ws.Range("K20").Value = '_res_.$Step2[4]'

$ws.Range("J21").Value = '_res_.$Step2[3]'
$ws.Range("K21").Value = '_res_.$Step2[4]'

# Row 22 holds the expected numeric results for that test - it previously
# only had a single (wrong) column; now it needs three numeric results.
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 3

# Leave the selection where the author ended up editing.
$ws.Range("K21").Select() | Out-Null
